# Update the cryptocurrency price/volume table to reflect the latest
# GitHub Actions data refresh.
#
# Values are written with a leading apostrophe so Excel treats them as
# literal text (matching the workbook's existing inlineStr cells) instead
# of auto-converting numeric-looking strings (e.g. "1.00", "0.581") into
# floating point numbers, which would silently drop formatting such as
# trailing zeros.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''98.636.08'
$ws.Range("E2").Value = '''  +0.63%  '
$ws.Range("D3").Value = '''3.477.78'
$ws.Range("E3").Value = '''  +2.28%  '
$ws.Range("D4").Value = '''1.00'
$ws.Range("E4").Value = '''  -0.03%  '
$ws.Range("D5").Value = '''256.08'
$ws.Range("E5").Value = '''  +0.83%  '
$ws.Range("D6").Value = '''674.87'
$ws.Range("E6").Value = '''  -0.38%  '
$ws.Range("E7").Value = '''  +4.69%  '
$ws.Range("D8").Value = '''0.436'
$ws.Range("E8").Value = '''  +1.41%  '
$ws.Range("E9").Value = '''  +1.86%  '
$ws.Range("E10").Value = '''  +0.03%  '
$ws.Range("D11").Value = '''3.476.22'
$ws.Range("E11").Value = '''  +2.32%  '
$ws.Range("D12").Value = '''46.56'
$ws.Range("E12").Value = '''  +11.92%  '
$ws.Range("E13").Value = '''  -0.99%  '
$ws.Range("D14").Value = '''98.499.90'
$ws.Range("E14").Value = '''  +0.79%  '
$ws.Range("D15").Value = '''6.23'
$ws.Range("E15").Value = '''  -0.47%  '
$ws.Range("E16").Value = '''  -0.92%  '
$ws.Range("D17").Value = '''4.130.95'
$ws.Range("E17").Value = '''  +2.57%  '
$ws.Range("D18").Value = '''9.13'
$ws.Range("E18").Value = '''  +2.42%  '
$ws.Range("D19").Value = '''3.473.50'
$ws.Range("E19").Value = '''  +1.70%  '
$ws.Range("D20").Value = '''18.60'
$ws.Range("E20").Value = '''  +7.13%  '
$ws.Range("D21").Value = '''0.540'
$ws.Range("E21").Value = '''  -4.55%  '
$ws.Range("D22").Value = '''11.85'
$ws.Range("E22").Value = '''  +7.42%  '
$ws.Range("D23").Value = '''524.87'
$ws.Range("E23").Value = '''  +3.49%  '
$ws.Range("D24").Value = '''3.47'
$ws.Range("E24").Value = '''  +0.93%  '
$ws.Range("D25").Value = '''0.0000205'
$ws.Range("E25").Value = '''  +0.85%  '
$ws.Range("D26").Value = '''6.91'
$ws.Range("E26").Value = '''  +5.66%  '
$ws.Range("D27").Value = '''98.63'
$ws.Range("E27").Value = '''  -0.61%  '
$ws.Range("D28").Value = '''12.78'
$ws.Range("E28").Value = '''  +1.66%  '
$ws.Range("D29").Value = '''3.644.38'
$ws.Range("E29").Value = '''  +1.78%  '
$ws.Range("D30").Value = '''3.14'
$ws.Range("E30").Value = '''  +21.18%  '
$ws.Range("D31").Value = '''12.36'
$ws.Range("E31").Value = '''  +7.36%  '
$ws.Range("E32").Value = '''  -2.28%  '
$ws.Range("E33").Value = '''  +0.02%  '
$ws.Range("D34").Value = '''0.193'
$ws.Range("E34").Value = '''  -0.58%  '
$ws.Range("B35").Value = '''PolygonEcosystemToken'
$ws.Range("C35").Value = '''https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Range("D35").Value = '''0.581'
$ws.Range("E35").Value = '''  +2.42%  '
$ws.Range("B36").Value = '''Binance-PegBSC-USD'
$ws.Range("C36").Value = '''https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range("D36").Value = '''1.01'
$ws.Range("E36").Value = '''  +1.11%  '
$ws.Range("D37").Value = '''30.35'
$ws.Range("E37").Value = '''  +3.22%  '
$ws.Range("D38").Value = '''8.23'
$ws.Range("E38").Value = '''  +4.15%  '
$ws.Range("E39").Value = '''  +1.80%  '
$ws.Range("D40").Value = '''538.68'
$ws.Range("E40").Value = '''  +1.20%  '
$ws.Range("D41").Value = '''0.157'
$ws.Range("E41").Value = '''  +2.88%  '
$ws.Range("E42").Value = '''  +0.02%  '
$ws.Range("D43").Value = '''0.894'
$ws.Range("E43").Value = '''  +2.97%  '
$ws.Range("D44").Value = '''1.80'
$ws.Range("E44").Value = '''  +4.16%  '
$ws.Range("E45").Value = '''  +2.13%  '
$ws.Range("E46").Value = '''  -1.11%  '
$ws.Range("D47").Value = '''5.81'
$ws.Range("E47").Value = '''  +1.38%  '
$ws.Range("D48").Value = '''8.81'
$ws.Range("E48").Value = '''  -1.46%  '
$ws.Range("D49").Value = '''3.67'
$ws.Range("E49").Value = '''  -2.12%  '
$ws.Range("E50").Value = '''  +7.73%  '
$ws.Range("D51").Value = '''56.12'
$ws.Range("E51").Value = '''  +0.96%  '
